# daily auto push: 2026-01-17 13:37 UTC
#
# The log table on Sheet1 gets a new row for 2026/01/17 (weekday 土) with
# time-bucket 19 and ranking 200, inserted directly above the existing
# "2026/12/29" block (row 669) -- this pushes every subsequent row down by
# one (old row 669 -> 670, ..., old row 710 -> 711).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 669:710 down to 670:711, leaving a blank row 669 to fill in.
$ws.Rows.Item(669).Insert()

# Column A/B of the new row duplicate the date/weekday already sitting in
# row 668 ("2026/01/17" / "土"). Copy them in (instead of re-typing the
# date string) so the date keeps being stored as literal text instead of
# being auto-parsed into a serial date value/format.
$ws.Range("A668:B668").Copy($ws.Range("A669:B669"))

# Numeric columns are plain numbers - safe to assign directly.
$ws.Cells.Item(669, 3).Value = 19
$ws.Cells.Item(669, 4).Value = 200
